$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = -12.214
$ws.Range("C12").Value = -12.53
$ws.Range("D13").Value = -7.831999999999999
$ws.Range("C18").Value = -12.283
